$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-mark price cells whose new values look numeric as Text so Excel
# keeps storing them as strings (matching the source inlineStr cells)
# instead of silently parsing them into floating point numbers.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '40.057.85'
$ws.Range('E2').Value = '  -3.89%  '
$ws.Range('D3').Value = '2.351.37'
$ws.Range('E3').Value = '  -5.08%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '309.34'
$ws.Range('E5').Value = '  -3.19%  '
$ws.Range('D6').Value = '84.23'
$ws.Range('E6').Value = '  -8.83%  '
$ws.Range('D7').Value = '0.530'
$ws.Range('E7').Value = '  -3.76%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').Value = '0.484'
$ws.Range('E9').Value = '  -5.35%  '
$ws.Range('D10').Value = '0.0822'
$ws.Range('E10').Value = '  -4.58%  '
$ws.Range('D11').Value = '30.07'
$ws.Range('E11').Value = '  -9.09%  '
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').Value = '2.717.68'
$ws.Range('E13').Value = '  -4.90%  '
$ws.Range('D14').Value = '6.41'
$ws.Range('E14').Value = '  -6.96%  '
$ws.Range('D15').Value = '14.85'
$ws.Range('E15').Value = '  -4.33%  '
$ws.Range('D16').Value = '2.367.44'
$ws.Range('E16').Value = '  -5.06%  '
$ws.Range('D17').Value = '0.751'
$ws.Range('E17').Value = '  -5.46%  '
$ws.Range('D18').Value = '40.149.70'
$ws.Range('E18').Value = '  -3.50%  '
$ws.Range('D19').Value = '0.0₃0902'
$ws.Range('E19').Value = '  -4.32%  '
$ws.Range('D20').Value = '6.08'
$ws.Range('E20').Value = '  -5.74%  '
$ws.Range('D21').Value = '67.95'
$ws.Range('E21').Value = '  -3.80%  '
$ws.Range('D22').Value = '10.70'
$ws.Range('E22').Value = '  -4.87%  '
$ws.Range('D23').Value = '233.78'
$ws.Range('E23').Value = '  -2.56%  '
$ws.Range('D24').Value = '2.56'
$ws.Range('E24').Value = '  -7.13%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('D26').Value = '1.79'
$ws.Range('E26').Value = '  -8.08%  '
$ws.Range('D27').Value = '23.49'
$ws.Range('E27').Value = '  -6.03%  '
$ws.Range('D28').Value = '2.21'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').Value = '9.21'
$ws.Range('E29').Value = '  -5.61%  '
$ws.Range('D30').Value = '34.29'
$ws.Range('E30').Value = '  -6.56%  '
$ws.Range('D31').Value = '152.30'
$ws.Range('E31').Value = '  -2.98%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('D33').Value = '5.17'
$ws.Range('E33').Value = '  -4.89%  '
$ws.Range('D34').Value = '0.0726'
$ws.Range('E34').Value = '  -5.04%  '
$ws.Range('E35').Value = '  -5.44%  '
$ws.Range('E36').Value = '  -2.52%  '
$ws.Range('D37').Value = '2.77'
$ws.Range('E37').Value = '  -4.13%  '
$ws.Range('D38').Value = '0.0993'
$ws.Range('E38').Value = '  -4.35%  '
$ws.Range('D39').Value = '15.71'
$ws.Range('E39').Value = '  -8.55%  '
$ws.Range('D40').Value = '1.70'
$ws.Range('E40').Value = '  -7.86%  '
$ws.Range('D41').Value = '3.81'
$ws.Range('E41').Value = '  -4.89%  '
$ws.Range('D42').Value = '2.37'
$ws.Range('E42').Value = '  -3.57%  '
$ws.Range('D43').Value = '1.960.11'
$ws.Range('E43').Value = '  -2.07%  '
$ws.Range('D44').Value = '0.0266'
$ws.Range('E44').Value = '  -6.38%  '
$ws.Range('D45').Value = '17.53'
$ws.Range('E45').Value = '  -6.36%  '
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').Value = '2.66'
$ws.Range('E47').Value = '  -10.83%  '
$ws.Range('D48').Value = '2.594.66'
$ws.Range('E48').Value = '  -4.46%  '
$ws.Range('D49').Value = '92.37'
$ws.Range('E49').Value = '  -5.50%  '
$ws.Range('D50').Value = '70.99'
$ws.Range('E50').Value = '  -6.14%  '
$ws.Range('D51').Value = '49.94'
$ws.Range('E51').Value = '  -4.70%  '

# Clear the temporary Text number format back to the default style so the
# cells keep no explicit style reference, same as in the original workbook.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
